$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet (sheet1): append Week 15 simulation results to the long
# space-separated result lists stored in B2/C2 (R - OFF/DEF) and B3/C3 (P - OFF/DEF)
# ---------------------------------------------------------------------------
$ydsSheet = $wb.Worksheets.Item("YDS")

$ydsSheet.Range("B2").Value = "5 3 5 3 3 5 2 5 1 7 0 6 6 1 8 2 7 3 2 10 3 8 10 4 2 8 4 6 0 3 1 3 1 5 -1 4 3 0 2 3 1 4 9 21 -1 11 3 4 3 3 6 3 4 6 8 0 1 7 7 3 8 6 1 4 6 4 1 -2 3 2 1 17 4 8 4 4 4 6 4 -1 13 13 3 6 1 4 5 5 4 6 2 4 0 9 4 3 6 0 2 3 2 2 20 5 6 2 7 0 1 9 2 2 0 1 -1 3 15 2 4 5 12 3 9 4 11 6 8 12 2 -1 1 12 19 12 9 5 3 1 2 -3 -1 4 13 9 -1 2 -1 0 24 6 2 5 6 3 5 6 0 3 5 8 2 1 3 4 8 6 7 10 7 6 3 2 5 10 3 3 4 6 2 7 -1 3 4 9 2 2 1 3 -3 2 -1 4 7 3 7 -3 4 2 2 1 11 6 12 4 1 4 5 8 8 1 2 3 1 11 8 3 28 3 4 2 2 6 6 5 4 3 1 3 4 9 -4 1 9 0 3 47 4 2 12 2 3 1 5 -1 8 2 1 1 2 9 1 -2 2 0 7 21 1 4 3 1"

$ydsSheet.Range("C2").Value = "4 2 4 0 2 7 1 0 3 -1 13 3 0 2 4 1 4 11 7 -1 0 2 3 10 0 0 8 2 0 8 1 -3 1 8 6 0 3 6 4 0 4 4 2 5 1 4 -2 0 2 11 2 3 1 4 1 5 15 0 6 -4 0 0 0 -1 1 2 4 0 3 1 13 7 4 5 -2 11 1 1 -1 5 6 2 5 -2 9 11 23 14 7 3 -1 2 0 29 1 1 6 4 0 6 13 5 12 6 0 6 0 8 8 1 2 -1 4 -1 -1 1 1 8 2 6 16 20 10 1 4 3 13 3 1 3 -1 3 3 11 7 1 1 3 2 15 2 4 2 3 6 5 1 22 1 3 1 3 1 5 -1 -8 10 5 2 3 15 4 4 -5 -1 -3 12 9 1 1 1 2 3 4 7 3 2 4 2 -1 5 0 3 2 8 6 0 1 3 16 1 2 12 0 9 6 0 4 -2 0 5 12 5 1 16 3 2 5 15 5 10 15 1 3 4 7 39 17 6 5 1 5 11 0 4 2 0 4 0 3 1 0 2 3 0 14 0 14 -3 10 22 1 29 18 4 11 5 18 23 2 0 7 8 3"

$ydsSheet.Range("B3").Value = "10 19 28 16 8 5 27 5 19 17 2 10 47 20 2 2 11 10 1 5 7 9 7 9 7 3 11 8 4 6 20 24 20 25 20 16 16 9 5 8 1 9 5 26 17 21 3 8 1 2 9 10 26 5 12 6 8 5 1 15 17 4 3 7 7 0 17 5 26 3 20 20 32 14 8 20 15 11 6 7 17 11 -4 23 13 31 17 8 -2 11 8 5 5 21 -1 9 9 7 8 6 13 13 28 3 6 0 27 8 1 4 28 26 11 11 9 7 23 18 11 8 9 26 13 10 10 18 6 15 62 6 10 10 11 12 3 4 11 5 12 6 14 12 16 4 15 7 34 19 9 22 9 22 16 10 2 8 8 11 3 9 23 -4 10 12 14 18 10 3 6 -4 5 4 4 9 15 0 8 5 4 8 10 2 27 10 -4 6 4 21 4 11 9 10 46 2 9 14 0 9 8 9 7 12 0 34 8 25 12 7 8 8 4 20 31 4 9 7 44 6 7 16 6 6 41 9 7 -1 8 2 12 8 50 11 2 3 29 9 8 3 6 7 11 4 5 14 15 8 6 2 15 4 0 4 3 40 10 8 15 3 13 7 10 8 6 6 6 7 4 35 7 12 16 3 13 9 17 14 14 18 4 5 8 3 21 5 3 20 9 0 0 4 -1 10 8 7 7 15 9 26 5 4 11 10 32 9 2 9 3 8 6 13 15 5 11 0 6 3 15 7 4 3 25 36 3 0 15 6 17 12 27 1 9 1 7 11 5 7 18 1 5 11 14 6 6 6 36 5 3 6 15 5 7 5 9 20 8 16 23 18 5 13 16 13 7 5 20 5 6 24 4 18 0 8 5 11 10 6 -7 14 58"

$ydsSheet.Range("C3").Value = "28 5 16 11 8 12 7 8 7 22 2 5 22 6 8 3 4 13 -5 4 11 7 10 5 16 6 14 6 4 3 7 10 9 6 21 2 3 6 11 9 31 10 0 9 5 9 4 5 23 24 4 8 24 7 6 1 14 8 3 8 14 4 23 12 8 7 22 5 7 8 1 2 4 10 2 6 7 6 12 22 11 5 6 6 6 9 2 11 21 22 2 75 0 20 5 40 10 5 17 6 4 10 3 7 4 5 7 17 1 16 15 5 16 7 3 11 4 12 15 9 4 4 21 13 5 10 12 15 1 12 3 21 30 0 9 7 7 21 4 34 9 12 10 5 24 5 3 12 5 2 21 6 10 8 4 23 10 2 7 1 0 11 5 21 23 5 6 4 16 2 1 6 -3 24 25 6 6 5 5 5 3 6 3 17 4 12 4 10 22 20 6 22 12 -6 11 8 5 6 4 5 18 7 16 5 14 12 14 -1 15 1 38 2 14 5 9 4 16 4 8 14 13 4 8 11 1 20 22 7 9 8 14 1 15 6 18 21 9 6 16 6 5 6 6 9 12 9 4 3 5 5 0 11 21 4 28 8 11 8 -1 7 2 4 5 5 8 4 9 8 62 10 14 3 3 15 9 5 12 22 9 23 -3 15 5 15 11 11 7 8 7 10 0 6 6 12 19 29 7 10 9 8 9 4 -3 14 2 19 6 9 4 16 6 20 6 17 8 8 8 2 7 30 1 10 11 5 3 25 24 21 9 5 6 10 12 6 1 13 9 -2 10 2 4 10 15 9 12 12 3 4 4 4 2 15 1 24 3"

# ---------------------------------------------------------------------------
# OFF sheet (sheet2): update the Home/Road totals with Week 15 numbers
# ---------------------------------------------------------------------------
$offSheet = $wb.Worksheets.Item("OFF")

$offSheet.Range("C2").Value = 150
$offSheet.Range("D2").Value = 9
$offSheet.Range("F2").Value = 16
$offSheet.Range("G2").Value = 36
$offSheet.Range("J2").Value = 27
$offSheet.Range("L2").Value = 258
$offSheet.Range("M2").Value = 170
$offSheet.Range("O2").Value = 13
$offSheet.Range("Q2").Value = 398

$offSheet.Range("C3").Value = 171
$offSheet.Range("D3").Value = 8
$offSheet.Range("E3").Value = 20
$offSheet.Range("F3").Value = 130
$offSheet.Range("G3").Value = 48
$offSheet.Range("H3").Value = 11
$offSheet.Range("I3").Value = 57
$offSheet.Range("J3").Value = 46
$offSheet.Range("N3").Value = 8

# ---------------------------------------------------------------------------
# DEF sheet (sheet4): update the Home/Road totals with Week 15 numbers
# ---------------------------------------------------------------------------
$defSheet = $wb.Worksheets.Item("DEF")

$defSheet.Range("C2").Value = 125
$defSheet.Range("D2").Value = 5
$defSheet.Range("E2").Value = 6
$defSheet.Range("F2").Value = 40
$defSheet.Range("G2").Value = 30
$defSheet.Range("J2").Value = 19
$defSheet.Range("L2").Value = 268
$defSheet.Range("M2").Value = 183
$defSheet.Range("O2").Value = 19
$defSheet.Range("P2").Value = 11
$defSheet.Range("Q2").Value = 393

$defSheet.Range("B3").Value = 4
$defSheet.Range("C3").Value = 192
$defSheet.Range("D3").Value = 8
$defSheet.Range("E3").Value = 27
$defSheet.Range("F3").Value = 104
$defSheet.Range("G3").Value = 40
$defSheet.Range("H3").Value = 20
$defSheet.Range("I3").Value = 45
$defSheet.Range("J3").Value = 57
$defSheet.Range("N3").Value = 14

# ---------------------------------------------------------------------------
# ST sheet (sheet6): numeric totals + the "D"/"RA"/"RM" stat lists
# ---------------------------------------------------------------------------
$stSheet = $wb.Worksheets.Item("ST")

$stSheet.Range("B2").Value = 80
$stSheet.Range("D2").Value = 43
$stSheet.Range("F2").Value = 455
$stSheet.Range("G2").Value = 437
$stSheet.Range("J2").Value = 183
$stSheet.Range("K2").Value = 173

$stSheet.Range("B3").Value = 66

$stSheet.Range("B4").Value = "51 50 44 56 76 64 67 74 73 64 64 69 67 69"
$stSheet.Range("D3").Value = "65 42 49 41 48 41 43 57 38 53 39 15 55 38 34 43 41 41 44 42 49 42 46 54 32 40 43 46 35 41 55 40 14 54 55 60 41 32 55 20 39 35 47"

$stSheet.Range("B5").Value = "13 19 14 45 27 0 25 33 36 28 19 72 30 29"
$stSheet.Range("D4").Value = "0 3 12 0 6 0 5 0 9 0 0 0 0 0 11 6 0 0 0 8 0 13 0 0 0 12 10 0 0 0 0 0 0 9 13 7 0 0 4 0 0 0 5"

$stSheet.Range("B6").Value = "20 41 31 16 24 19 24 24 13 10 17 17 17 19 34 18 13 15 24 30 12 15 33 7"
$stSheet.Range("D5").Value = "0 0 10 0 0 0 0 0"

# ---------------------------------------------------------------------------
# TURNS sheet (sheet7): Home/Road turnover totals
# ---------------------------------------------------------------------------
$turnsSheet = $wb.Worksheets.Item("TURNS")

$turnsSheet.Range("C2").Value = 11
$turnsSheet.Range("D2").Value = 6
$turnsSheet.Range("E3").Value = 8

# ---------------------------------------------------------------------------
# Move the active tab from ST back to YDS (first sheet)
# ---------------------------------------------------------------------------
$ydsSheet.Activate()
